$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: becomes the "Tretåig hackspett" observation (previously row 12's species data),
# with its own Id/Taxonsorteringsordning/coordinates, and gains a "Ringhack" comment.
$ws.Range("A11").Value2 = 111221736
$ws.Range("B11").Value2 = 56398
$ws.Range("E11").Value2 = 100109
$ws.Range("F11").Value2 = "Tretåig hackspett"
$ws.Range("G11").Value2 = "Picoides tridactylus"
$ws.Range("H11").Value2 = "(Linnaeus, 1758)"
$ws.Range("I11").ClearContents()
$ws.Range("M11").ClearContents()
$ws.Range("Q11").Value2 = 443249.6264723797
$ws.Range("R11").Value2 = 6909840.911127058
$ws.Range("AC11").Value2 = "Ringhack"

# Row 12: keeps its species data, but its Id and coordinates are updated.
$ws.Range("A12").Value2 = 111221709
$ws.Range("Q12").Value2 = 443254.9775056695
$ws.Range("R12").Value2 = 6909826.869210822

# Row 13: becomes the "Talltita" observation (previously row 11's species data),
# with its own Id/Taxonsorteringsordning/coordinates, and loses the "Ringhack" comment.
$ws.Range("A13").Value2 = 111221699
$ws.Range("B13").Value2 = 56543
$ws.Range("E13").Value2 = 103021
$ws.Range("F13").Value2 = "Talltita"
$ws.Range("G13").Value2 = "Poecile montanus"
$ws.Range("H13").Value2 = "(Conrad von Baldenstein, 1827)"

# "1" must be stored as text (like the original data), not a number, so build it
# as a text formula result and freeze it to a value to avoid numeric auto-conversion.
$ws.Range("I13").Formula = "=""1"""
$ws.Range("I13").Copy() | Out-Null
$ws.Range("I13").PasteSpecial(-4163) | Out-Null   # xlPasteValues
$excel.CutCopyMode = $false

$ws.Range("M13").Value2 = "spel/sång"
$ws.Range("Q13").Value2 = 443097.6233577073
$ws.Range("R13").Value2 = 6909995.088246249
$ws.Range("AC13").ClearContents()
